{"js": "// Remove the last column (the \"RECOVERY\" column) from the table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst columnCount = table.values[0].length;\ntable.deleteColumns(columnCount - 1, 1);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$t.Columns($t.Columns.Count).Delete()\n"}
